$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 591622.75
$ws.Range("J103").Value = 929342.9
$ws.Range("L103").Value = 2788028.7
$ws.Range("N103").Value = -2789200.7
$ws.Range("H112").Value = 875
$ws.Range("J112").Value = 896.9048
$ws.Range("L112").Value = 2690.7144
$ws.Range("N112").Value = -4906.7144
$ws.Range("H116").Value = 4946.8438
$ws.Range("J116").Value = 1864.7
$ws.Range("L116").Value = 1864.7
$ws.Range("N116").Value = -8748.700000000001
$ws.Range("H137").Value = 1622.6666
$ws.Range("I137").Value = 1096.5385
$ws.Range("J137").Value = 1964.65
$ws.Range("K137").Value = 3289.6155
$ws.Range("L137").Value = 5893.950000000001
$ws.Range("M137").Value = -739.6155000000003
$ws.Range("N137").Value = -10993.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3617.3086
$ws.Range("I32").Value = 3006.6716
$ws.Range("J32").Value = 6539.643
$ws.Range("K32").Value = 3006.6716
$ws.Range("L32").Value = 6539.643
$ws.Range("M32").Value = -2719.6716
$ws.Range("N32").Value = -7113.643
$ws.Range("H74").Value = 2111.9565
$ws.Range("I74").Value = 2073.3
$ws.Range("J74").Value = 2141.6924
$ws.Range("K74").Value = 2073.3
$ws.Range("L74").Value = 2141.6924
$ws.Range("M74").Value = -1199.3
$ws.Range("N74").Value = -3889.6924
$ws.Range("H77").Value = 2111.9565
$ws.Range("I77").Value = 2073.3
$ws.Range("J77").Value = 2141.6924
$ws.Range("K77").Value = 10366.5
$ws.Range("L77").Value = 10708.462
$ws.Range("M77").Value = -5998.5
$ws.Range("N77").Value = -19444.462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2360.2236
$ws.Range("I31").Value = 1512.2778
$ws.Range("J31").Value = 2623.3794
$ws.Range("K31").Value = 1512.2778
$ws.Range("L31").Value = 2623.3794
$ws.Range("M31").Value = -1217.2778
$ws.Range("N31").Value = -3213.3794
$ws.Range("H34").Value = 2360.2236
$ws.Range("I34").Value = 1512.2778
$ws.Range("J34").Value = 2623.3794
$ws.Range("K34").Value = 1512.2778
$ws.Range("L34").Value = 2623.3794
$ws.Range("M34").Value = -1310.2778
$ws.Range("N34").Value = -3027.3794
$ws.Range("H134").Value = 3824.158
$ws.Range("I134").Value = 3814.389
$ws.Range("K134").Value = 11443.167
$ws.Range("M134").Value = -8908.167000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 139049.8
$ws.Range("I5").Value = 11459.333
$ws.Range("J5").Value = 171858.77
$ws.Range("K5").Value = 34377.999
$ws.Range("L5").Value = 515576.3099999999
$ws.Range("M5").Value = -34265.999
$ws.Range("N5").Value = -515800.3099999999
$ws.Range("H14").Value = 87.45
$ws.Range("I14").Value = 87.45
$ws.Range("K14").Value = 262.35
$ws.Range("M14").Value = -89.35000000000002
$ws.Range("H33").Value = 10030048
$ws.Range("I33").Value = 50000024
$ws.Range("J33").Value = 37553.75
$ws.Range("K33").Value = 300000144
$ws.Range("L33").Value = 225322.5
$ws.Range("M33").Value = -299999861
$ws.Range("N33").Value = -225888.5
$ws.Range("H97").Value = 10000396
$ws.Range("J97").Value = 486.33334
$ws.Range("L97").Value = 1459.00002
$ws.Range("N97").Value = -2451.00002
$ws.Range("H98").Value = 7143164.5
$ws.Range("I98").Value = 350
$ws.Range("J98").Value = 8333633.5
$ws.Range("K98").Value = 1050
$ws.Range("L98").Value = 25000900.5
$ws.Range("M98").Value = 448
$ws.Range("N98").Value = -25003896.5
$ws.Range("H102").Value = 4800
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 4800
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 14400
$ws.Range("M102").Value = ""
$ws.Range("N102").Value = -19268
$ws.Range("H117").Value = 15879760
$ws.Range("I117").Value = 12957.25
$ws.Range("J117").Value = 25643946
$ws.Range("K117").Value = 38871.75
$ws.Range("L117").Value = 76931838
$ws.Range("M117").Value = -35429.75
$ws.Range("N117").Value = -76938722
$ws.Range("H121").Value = 899.36664
$ws.Range("J121").Value = 1005.7451
$ws.Range("L121").Value = 3017.2353
$ws.Range("N121").Value = -5637.2353
$ws.Range("H129").Value = 1665.5312
$ws.Range("I129").Value = 879.93335
$ws.Range("J129").Value = 2358.7058
$ws.Range("K129").Value = 2639.80005
$ws.Range("L129").Value = 7076.117400000001
$ws.Range("M129").Value = 2360.19995
$ws.Range("N129").Value = -17076.1174
$ws.Range("H135").Value = 139049.8
$ws.Range("I135").Value = 11459.333
$ws.Range("J135").Value = 171858.77
$ws.Range("K135").Value = 103133.997
$ws.Range("L135").Value = 1546728.93
$ws.Range("M135").Value = -100598.997
$ws.Range("N135").Value = -1551798.93

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3474603.5
$ws.Range("I22").Value = 12346868
$ws.Range("K22").Value = 12346868
$ws.Range("M22").Value = -12346573
$ws.Range("H27").Value = 3474603.5
$ws.Range("I27").Value = 12346868
$ws.Range("K27").Value = 12346868
$ws.Range("M27").Value = -12346761
$ws.Range("H46").Value = 17545342
$ws.Range("I46").Value = 27778782
$ws.Range("J46").Value = 2299.8572
$ws.Range("K46").Value = 27778782
$ws.Range("L46").Value = 2299.8572
$ws.Range("M46").Value = -27778594
$ws.Range("N46").Value = -2675.8572
$ws.Range("H55").Value = 17857490
$ws.Range("I55").Value = 243.5
$ws.Range("J55").Value = 35714736
$ws.Range("K55").Value = 243.5
$ws.Range("L55").Value = 35714736
$ws.Range("M55").Value = -70.5
$ws.Range("N55").Value = -35715082
$ws.Range("H82").Value = 7417607.5
$ws.Range("I82").Value = 1518.4445
$ws.Range("J82").Value = 15760708
$ws.Range("K82").Value = 1518.4445
$ws.Range("L82").Value = 15760708
$ws.Range("M82").Value = -1157.4445
$ws.Range("N82").Value = -15761430
$ws.Range("H85").Value = 7417607.5
$ws.Range("I85").Value = 1518.4445
$ws.Range("J85").Value = 15760708
$ws.Range("K85").Value = 1518.4445
$ws.Range("L85").Value = 15760708
$ws.Range("M85").Value = -270.4445000000001
$ws.Range("N85").Value = -15763204
$ws.Range("H100").Value = 2639.6155
$ws.Range("I100").Value = 1901.25
$ws.Range("J100").Value = 2967.7778
$ws.Range("K100").Value = 1901.25
$ws.Range("L100").Value = 2967.7778
$ws.Range("M100").Value = -1360.25
$ws.Range("N100").Value = -4049.7778
$ws.Range("H132").Value = 21670518
$ws.Range("I132").Value = 25494326
$ws.Range("J132").Value = 2268.3333
$ws.Range("K132").Value = 76482978
$ws.Range("L132").Value = 6804.999899999999
$ws.Range("M132").Value = -76480448
$ws.Range("N132").Value = -11864.9999
$ws.Range("H136").Value = 3730.6553
$ws.Range("I136").Value = 2289.7273
$ws.Range("J136").Value = 8259.286
$ws.Range("K136").Value = 6869.1819
$ws.Range("L136").Value = 24777.858
$ws.Range("M136").Value = -4319.1819
$ws.Range("N136").Value = -29877.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = ""
$ws.Range("N26").Value = ""
